# Apply the commit's edits to the "2017 MCM Problem C" workbook.
#
# Summary of the change (per the OOXML diff):
#   - Sheet "parsed mile posts": the stray value in I2 (100.93) is removed
#     and a brand new data row (row 3) is appended with the same shape as
#     row 2 (columns A-G), reusing the existing shared string "IS" for E3.
#   - The used range / dimension grows from A1:I2 to A1:I3.
#   - The selection (view state) on both sheets is updated to reflect where
#     the user was working (A4:J5 on "parsed mile posts", and a combination
#     of A4:J5 / B12 on "definitions").

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("parsed mile posts")
$wsDefs = $wb.Worksheets.Item("definitions")

# --- Update sheet data -----------------------------------------------------

# Remove the now-orphaned value in I2 (was 100.93).
$wsData.Range("I2").ClearContents()

# Append the new row of observations as row 3 (mirrors row 2's layout).
$wsData.Range("A3").Value = 5
$wsData.Range("B3").Value = 10.56
$wsData.Range("C3").Value = 10.93
$wsData.Range("D3").Value = 177000
$wsData.Range("E3").Value = "IS"
$wsData.Range("F3").Value = 2
$wsData.Range("G3").Value = 2

# --- Update view/selection state --------------------------------------------

$null = $wsData.Activate()
$null = $wsData.Range("A4:J5").Select()

$null = $wsDefs.Activate()
$null = $wsDefs.Range("B12").Select()

# Leave focus back on the data sheet, matching tabSelected="true" in the
# original workbook.
$null = $wsData.Activate()
